# Mark the completed rubric rows with an "X" in column F, matching the
# grader's workflow of checking off additional custom-feature rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = "X"
$ws.Range("F29").Value = "X"
$ws.Range("F30").Value = "X"
$ws.Range("F31").Value = "X"
$ws.Range("F74").Value = "X"

# Leave the selection where the grader last clicked.
$ws.Range("F31").Select() | Out-Null
